$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "descripcion" column (column I). Excel shifts the
# following columns (categoria, marca, clase, ubicacion, proveedor) left
# by one, updates the used range / column widths automatically.
$ws.Columns("I:I").Delete()

# Move the active selection to M2 (matches the post-edit cursor position).
$ws.Range("M2").Select() | Out-Null
